$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7791
$ws1.Range("F5").Value = 7791
$ws1.Range("F8").Value = 2114
$ws1.Range("F9").Value = 8586
$ws1.Range("F10").Value = 8586
$ws1.Range("F14").Value = 5736
$ws1.Range("F16").Value = 2697
$ws1.Range("F19").Value = 412
$ws1.Range("F21").Value = 40
$ws1.Range("F22").Value = 587
$ws1.Range("F23").Value = 14
$ws1.Range("F24").Value = 3766
$ws1.Range("F26").Value = 52
$ws1.Range("F28").Value = 14
$ws1.Range("F29").Value = 91
$ws1.Range("F30").Value = 7
$ws1.Range("F31").Value = 4092
$ws1.Range("F35").Value = 379
$ws1.Range("F36").Value = 151
$ws1.Range("F38").Value = 1534
$ws1.Range("F39").Value = 691
$ws1.Range("F42").Value = 3395
$ws1.Range("F43").Value = 62
$ws1.Range("F46").Value = 3347
$ws1.Range("F48").Value = 2317

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 64
$ws2.Range("F6").Value = 8

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7791
$ws4.Range("F6").Value = 7791
$ws4.Range("F8").Value = 2114
$ws4.Range("F9").Value = 8586
$ws4.Range("F10").Value = 8586
$ws4.Range("F13").Value = 5736
$ws4.Range("F15").Value = 2697
$ws4.Range("F18").Value = 412
$ws4.Range("F21").Value = 40
$ws4.Range("F23").Value = 587
$ws4.Range("F25").Value = 3766
$ws4.Range("F27").Value = 52
$ws4.Range("F29").Value = 14
$ws4.Range("F30").Value = 7
$ws4.Range("F31").Value = 4092
$ws4.Range("F34").Value = 379
$ws4.Range("F35").Value = 151
$ws4.Range("F37").Value = 64
$ws4.Range("F38").Value = 1534
$ws4.Range("F39").Value = 691
$ws4.Range("F43").Value = 3395
$ws4.Range("F44").Value = 62
$ws4.Range("F47").Value = 3347
$ws4.Range("F48").Value = 2317
